# "Generate Report for Handback"
#
# Re-running the handback-report generator picked up a newer handoff/handback
# round-trip for the second tracked file (79b83687-3609-4f18-b3c8-7a982e7f2368.*)
# and refreshed its timestamps:
#   - Overview!G3            "Latest HO Xliff Generate Date" rollup
#   - zh-cn!H3 / zh-cn!K3     "Correspond Handoff/Handback Datetime" for zh-cn
#   - de-de!H3 / de-de!K3     "Correspond Handoff/Handback Datetime" for de-de
#
# The first tracked file's row (27a7ff87-1504-46fb-a672-5bd3229e647d.*) and
# every other column/row/sheet is untouched.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-10-18 12:28:13"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-10-18 12:28:02"
$wsZhCn.Range("K3").Value = "2016-10-18 12:28:44"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-10-18 12:28:13"
$wsDeDe.Range("K3").Value = "2016-10-18 12:29:01"
